# Auto-update: dice_jobs_list.xlsx
# Append new job listing rows (61-68) to the dice jobs list worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Software Engineer III Backend Systems (AI & Data Platforms) | Redmond, WA (Onsite) | Contract |", "https://www.dice.com/job-detail/f2a13a78-2b8f-4c2d-a13f-d6edbe09fad8", "Redmond, Washington", "Contract", "Depends on Experience", "Tri-Force Consulting Services Inc"),
    @("SENIOR DEVOPS LEAD WITH Google Cloud Platform, KUBERNETES, TERRAFORM", "https://www.dice.com/job-detail/f4b94690-2253-4ba4-aec3-c6f1ea377eb2", "Remote", "Contract", "Depends on Experience", "Myticas LLC"),
    @("Lead Java AWS Developer with AI exp : Chicago, IL (In-Person interview)", "https://www.dice.com/job-detail/1114ad54-a1b1-4931-8891-755ee17a6c46", "Chicago, Illinois", "Contract, Third Party", "Depends on Experience", "Randstad Digital"),
    @("Microservices Sr Backend Developer Lead", "https://www.dice.com/job-detail/1ba91949-8f43-4ec0-bdfd-030ebe183bbb", "Pleasanton, California", "Contract, Third Party", "Depends on Experience", "Digitek Software, Inc."),
    @("Java Backend Engineer", "https://www.dice.com/job-detail/814659f0-f99d-413c-ab52-901914b6180c", "Hybrid in Dayton, Ohio", "Contract", "Depends on Experience", "Apetan Consulting"),
    @("SAP Vistex Consultant  |Chicago, IL (Hybrid) | Contract W2 only", "https://www.dice.com/job-detail/7a2c89a4-feb1-44dd-abbd-f8f9dc06a02c", "Hybrid in Chicago, Illinois", "Contract", "`$65 - `$70", "Robert Half"),
    @("Google Cloud Platform Cloud Engineer/Architect", "https://www.dice.com/job-detail/0ecc79bc-5620-41ac-9513-6201008d633e", "Remote", "Contract", "Depends on Experience", "Zen & Art"),
    @("Senior Cloud Engineer - Google Cloud Platform", "https://www.dice.com/job-detail/e3a7d768-9878-4f4c-a8b7-a443716d8640", "Remote", "Contract", "Up to `$80", "Sysazzle")
)

$startRow = 61
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
